$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.011815278060212
$ws.Range("D2").Value = 1.04523878805012
$ws.Range("E2").Value = 1.013940126994007
$ws.Range("F2").Value = 1.046357259005482
$ws.Range("I2").Value = 1.035314664144424
$ws.Range("J2").Value = 1.017061578419777
$ws.Range("K2").Value = 1.048007587784406
$ws.Range("L2").Value = 1.016799772201894
$ws.Range("M2").Value = 1.049122921490067
$ws.Range("N2").Value = 1.009802614112428

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.012733166602526
$ws.Range("D3").Value = 1.045849588508455
$ws.Range("E3").Value = 1.0147176133934
$ws.Range("F3").Value = 1.047181870607714
$ws.Range("I3").Value = 1.035405190088722
$ws.Range("J3").Value = 1.017613058064312
$ws.Range("K3").Value = 1.048430311104754
$ws.Range("L3").Value = 1.017382173085179
$ws.Range("M3").Value = 1.049759125557389
$ws.Range("N3").Value = 1.009988578144779

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.013327919343449
$ws.Range("D4").Value = 1.04624451286278
$ws.Range("E4").Value = 1.015221789810603
$ws.Range("F4").Value = 1.047715732527017
$ws.Range("I4").Value = 1.035462276868634
$ws.Range("J4").Value = 1.017970137795792
$ws.Range("K4").Value = 1.048702843343768
$ws.Range("L4").Value = 1.01775946164397
$ws.Range("M4").Value = 1.050170425969821
$ws.Range("N4").Value = 1.010108887449848

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.013578147984596
$ws.Range("D5").Value = 1.046410464012493
$ws.Range("E5").Value = 1.015434005204212
$ws.Range("F5").Value = 1.047940233618147
$ws.Range("I5").Value = 1.035485918733868
$ws.Range("J5").Value = 1.018120309473828
$ws.Range("K5").Value = 1.048817175231137
$ws.Range("L5").Value = 1.017918177062038
$ws.Range("M5").Value = 1.050343247033754
$ws.Range("N5").Value = 1.01015945982432

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.01362017385296
$ws.Range("D6").Value = 1.046438323476359
$ws.Range("E6").Value = 1.015469652283856
$ws.Range("F6").Value = 1.047977932117646
$ws.Range("I6").Value = 1.035489867324269
$ws.Range("J6").Value = 1.018145527175832
$ws.Range("K6").Value = 1.048836357894672
$ws.Range("L6").Value = 1.017944832104327
$ws.Range("M6").Value = 1.050372259137558
$ws.Range("N6").Value = 1.010167950804441

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.013331262147585
$ws.Range("D7").Value = 1.046246730607819
$ws.Range("E7").Value = 1.015224624424795
$ws.Range("F7").Value = 1.047718732066877
$ws.Range("I7").Value = 1.035462594177825
$ws.Range("J7").Value = 1.017972144180867
$ws.Range("K7").Value = 1.048704371999901
$ws.Range("L7").Value = 1.017761582002641
$ws.Range("M7").Value = 1.050172735568015
$ws.Range("N7").Value = 1.010109563223152

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.012125312853472
$ws.Range("D8").Value = 1.045445272642182
$ws.Range("E8").Value = 1.014202655421618
$ws.Range("F8").Value = 1.046635880407035
$ws.Range("I8").Value = 1.035345565745736
$ws.Range("J8").Value = 1.017247904163128
$ws.Range("K8").Value = 1.04815065471166
$ws.Range("L8").Value = 1.016996506300653
$ws.Range("M8").Value = 1.049338004266309
$ws.Range("N8").Value = 1.00986546579872

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.010006586513363
$ws.Range("D9").Value = 1.044030766115486
$ws.Range("E9").Value = 1.01241023523751
$ws.Range("F9").Value = 1.044730023750525
$ws.Range("I9").Value = 1.03512797431812
$ws.Range("J9").Value = 1.015973549956763
$ws.Range("K9").Value = 1.047167367144599
$ws.Range("L9").Value = 1.015651736225345
$ws.Range("M9").Value = 1.047864383613401
$ws.Range("N9").Value = 1.009435186291296

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.00859841104739
$ws.Range("D10").Value = 1.043086388388752
$ws.Range("E10").Value = 1.011221038377077
$ws.Range("F10").Value = 1.04346110878446
$ws.Range("I10").Value = 1.03497531549122
$ws.Range("J10").Value = 1.015125284726642
$ws.Range("K10").Value = 1.046506857738616
$ws.Range("L10").Value = 1.014757569936434
$ws.Range("M10").Value = 1.046880259800232
$ws.Range("N10").Value = 1.009148258544228

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.00798969117841
$ws.Range("D11").Value = 1.042677163960952
$ws.Range("E11").Value = 1.010707485525935
$ws.Range("F11").Value = 1.042912075287348
$ws.Range("I11").Value = 1.034907420557811
$ws.Range("J11").Value = 1.014758298615533
$ws.Range("K11").Value = 1.04621968938736
$ws.Range("L11").Value = 1.014370956847749
$ws.Range("M11").Value = 1.046453742357206
$ws.Range("N11").Value = 1.009024003223204

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.007763741242319
$ws.Range("D12").Value = 1.042525116519933
$ws.Range("E12").Value = 1.010516937554775
$ws.Range("F12").Value = 1.04270820440258
$ws.Range("I12").Value = 1.034881932827259
$ws.Range("J12").Value = 1.014622032473325
$ws.Range("K12").Value = 1.046112849224301
$ws.Range("L12").Value = 1.014227437817296
$ws.Range("M12").Value = 1.046295258871707
$ws.Range("N12").Value = 1.008977847625811

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.007812201206402
$ws.Range("D13").Value = 1.042557733134877
$ws.Range("E13").Value = 1.010557801285568
$ws.Range("F13").Value = 1.042751932440065
$ws.Range("I13").Value = 1.034887412170724
$ws.Range("J13").Value = 1.014651259797008
$ws.Range("K13").Value = 1.046135774605017
$ws.Range("L13").Value = 1.014258219213845
$ws.Range("M13").Value = 1.046329256615985
$ws.Range("N13").Value = 1.008987748225516

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.007971010903217
$ws.Range("D14").Value = 1.042664596542
$ws.Range("E14").Value = 1.010691730510686
$ws.Range("F14").Value = 1.042895221925598
$ws.Range("I14").Value = 1.034905319207367
$ws.Range("J14").Value = 1.014747033811528
$ws.Range("K14").Value = 1.046210861466308
$ws.Range("L14").Value = 1.014359091753218
$ws.Range("M14").Value = 1.046440643192124
$ws.Range("N14").Value = 1.009020188018116

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.008068879466227
$ws.Range("D15").Value = 1.042730432946341
$ws.Range("E15").Value = 1.010774276389823
$ws.Range("F15").Value = 1.042983515934193
$ws.Range("I15").Value = 1.034916316764717
$ws.Range("J15").Value = 1.014806049844883
$ws.Range("K15").Value = 1.046257102080689
$ws.Range("L15").Value = 1.014421254121224
$ws.Range("M15").Value = 1.046509264787458
$ws.Range("N15").Value = 1.009040175040454

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.008638831559846
$ws.Range("D16").Value = 1.043113541078001
$ws.Range("E16").Value = 1.011255150334092
$ws.Range("F16").Value = 1.0434975552659
$ws.Range("I16").Value = 1.034979783753235
$ws.Range("J16").Value = 1.01514964717957
$ws.Range("K16").Value = 1.046525891793011
$ws.Range("L16").Value = 1.014783240212172
$ws.Range("M16").Value = 1.046908558408389
$ws.Range("N16").Value = 1.009156504707719

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.008996623803698
$ws.Range("D17").Value = 1.043353775217882
$ws.Range("E17").Value = 1.011557159645651
$ws.Range("F17").Value = 1.043820111266133
$ws.Range("I17").Value = 1.035019115537021
$ws.Range("J17").Value = 1.015365262651369
$ws.Range("K17").Value = 1.046694186227457
$ws.Range("L17").Value = 1.015010457071249
$ws.Range("M17").Value = 1.047158922978259
$ws.Range("N17").Value = 1.009229471860757

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.009205417216198
$ws.Range("D18").Value = 1.043493870379669
$ws.Range("E18").Value = 1.011733449349629
$ws.Range("F18").Value = 1.044008292620518
$ws.Range("I18").Value = 1.035041884082608
$ws.Range("J18").Value = 1.015491058154398
$ws.Range("K18").Value = 1.046792237154242
$ws.Range("L18").Value = 1.015143043434719
$ws.Range("M18").Value = 1.047304919084211
$ws.Range("N18").Value = 1.009272031015224

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.009276627176207
$ws.Range("D19").Value = 1.04354163416455
$ws.Range("E19").Value = 1.011793582059192
$ws.Range("F19").Value = 1.044072464298159
$ws.Range("I19").Value = 1.035049618197067
$ws.Range("J19").Value = 1.015533956337404
$ws.Range("K19").Value = 1.046825650873602
$ws.Range("L19").Value = 1.015188261187497
$ws.Range("M19").Value = 1.047354693594556
$ws.Range("N19").Value = 1.009286542345028

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.008958225808403
$ws.Range("D20").Value = 1.043328003387879
$ws.Range("E20").Value = 1.01152474315953
$ws.Range("F20").Value = 1.043785499917028
$ws.Range("I20").Value = 1.035014913499204
$ws.Range("J20").Value = 1.015342125974597
$ws.Range("K20").Value = 1.046676141437694
$ws.Range("L20").Value = 1.014986073195528
$ws.Range("M20").Value = 1.04713206507439
$ws.Range("N20").Value = 1.009221643312527

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.007924241098285
$ws.Range("D21").Value = 1.042633129108873
$ws.Range("E21").Value = 1.010652285926348
$ws.Range("F21").Value = 1.042853024942179
$ws.Range("I21").Value = 1.03490005343773
$ws.Range("J21").Value = 1.014718829394995
$ws.Range("K21").Value = 1.046188755020268
$ws.Range("L21").Value = 1.014329384913619
$ws.Range("M21").Value = 1.046407844149408
$ws.Range("N21").Value = 1.009010635349586

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.007275035392883
$ws.Range("D22").Value = 1.042195984819542
$ws.Range("E22").Value = 1.010104944058872
$ws.Range("F22").Value = 1.042267115789236
$ws.Range("I22").Value = 1.034826283123856
$ws.Range("J22").Value = 1.014327221170403
$ws.Range("K22").Value = 1.045881315471179
$ws.Range("L22").Value = 1.013916998701952
$ws.Range("M22").Value = 1.045952175169117
$ws.Range("N22").Value = 1.008877956996187

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.007619105807618
$ws.Range("D23").Value = 1.042427746171739
$ws.Range("E23").Value = 1.010394985460836
$ws.Range("F23").Value = 1.04257768103514
$ws.Range("I23").Value = 1.034865537126152
$ws.Range("J23").Value = 1.014534792928617
$ws.Range("K23").Value = 1.046044389236752
$ws.Range("L23").Value = 1.014135564685808
$ws.Range("M23").Value = 1.046193763831568
$ws.Range("N23").Value = 1.008948293039783

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.008975575903081
$ws.Range("D24").Value = 1.043339648659459
$ws.Range("E24").Value = 1.011539390364161
$ws.Range("F24").Value = 1.043801139172838
$ws.Range("I24").Value = 1.035016812753904
$ws.Range("J24").Value = 1.015352580348339
$ws.Range("K24").Value = 1.046684295448856
$ws.Range("L24").Value = 1.014997091050753
$ws.Range("M24").Value = 1.047144201120568
$ws.Range("N24").Value = 1.00922518070072

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.010553573466338
$ws.Range("D25").Value = 1.044396701474052
$ws.Range("E25").Value = 1.012872611869741
$ws.Range("F25").Value = 1.045222451150467
$ws.Range("I25").Value = 1.035185569622356
$ws.Range("J25").Value = 1.016302775650487
$ws.Range("K25").Value = 1.047422456615328
$ws.Range("L25").Value = 1.015998982752112
$ws.Range("M25").Value = 1.048245659832574
$ws.Range("N25").Value = 1.009546438877705

